# Auto-generated edit script: updates FFXIV leve-profit price/profit columns
# (H/I/J/K/L/M/N) across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR per the scheduled
# price-refresh commit. Values below are the refreshed market-board prices
# and their downstream profit recalculations.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H11").Value = 44454.65
$ws_ALC.Range("I11").Value = 44454.65
$ws_ALC.Range("K11").Value = 44454.65
$ws_ALC.Range("M11").Value = -44314.65
$ws_ALC.Range("H62").Value = 8535.625
$ws_ALC.Range("I62").Value = 10448.5
$ws_ALC.Range("K62").Value = 10448.5
$ws_ALC.Range("M62").Value = -9824.5
$ws_ALC.Range("H65").Value = 8535.625
$ws_ALC.Range("I65").Value = 10448.5
$ws_ALC.Range("K65").Value = 52242.5
$ws_ALC.Range("M65").Value = -49122.5
$ws_ALC.Range("H96").Value = 1406.875
$ws_ALC.Range("I96").Value = 314.75
$ws_ALC.Range("J96").Value = 2499
$ws_ALC.Range("K96").Value = 944.25
$ws_ALC.Range("L96").Value = 7497
$ws_ALC.Range("M96").Value = 428.75
$ws_ALC.Range("N96").Value = -10243
$ws_ALC.Range("H131").Value = 31666.666
$ws_ALC.Range("I131").Value = 2500
$ws_ALC.Range("J131").Value = 90000
$ws_ALC.Range("K131").Value = 7500
$ws_ALC.Range("L131").Value = 270000
$ws_ALC.Range("M131").Value = -2460
$ws_ALC.Range("N131").Value = -280080
$ws_ALC.Range("H138").Value = 4420.758
$ws_ALC.Range("I138").Value = 4009.8
$ws_ALC.Range("J138").Value = 4494.143
$ws_ALC.Range("K138").Value = 12029.4
$ws_ALC.Range("L138").Value = 13482.429
$ws_ALC.Range("M138").Value = -6889.400000000001
$ws_ALC.Range("N138").Value = -23762.429
$ws_ALC.Range("H141").Value = 12568.4
$ws_ALC.Range("I141").Value = 11526.286
$ws_ALC.Range("K141").Value = 34578.858
$ws_ALC.Range("M141").Value = -29398.858

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H61").Value = 1452725.4
$ws_ARM.Range("J61").Value = 6253352.5
$ws_ARM.Range("L61").Value = 6253352.5
$ws_ARM.Range("N61").Value = -6253776.5
$ws_ARM.Range("H132").Value = 896892.9
$ws_ARM.Range("I132").Value = 1045832.8
$ws_ARM.Range("K132").Value = 3137498.4
$ws_ARM.Range("M132").Value = -3134968.4
$ws_ARM.Range("H133").Value = 94563.25
$ws_ARM.Range("J133").Value = 94563.25
$ws_ARM.Range("L133").Value = 94563.25
$ws_ARM.Range("N133").Value = -99623.25
$ws_ARM.Range("H136").Value = 1452725.4
$ws_ARM.Range("J136").Value = 6253352.5
$ws_ARM.Range("L136").Value = 18760057.5
$ws_ARM.Range("N136").Value = -18765157.5

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H35").Value = 36359.4
$ws_BSM.Range("J35").Value = 36359.4
$ws_BSM.Range("L35").Value = 36359.4
$ws_BSM.Range("N35").Value = -36979.4
$ws_BSM.Range("H105").Value = 6475.077
$ws_BSM.Range("I105").Value = 5512
$ws_BSM.Range("K105").Value = 5512
$ws_BSM.Range("M105").Value = -3765

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H16").Value = 79362.92
$ws_CRP.Range("I16").Value = 2374.2222
$ws_CRP.Range("K16").Value = 2374.2222
$ws_CRP.Range("M16").Value = -2087.2222
$ws_CRP.Range("H31").Value = 1738637
$ws_CRP.Range("I31").Value = 1918081.8
$ws_CRP.Range("J31").Value = 4004.6667
$ws_CRP.Range("K31").Value = 1918081.8
$ws_CRP.Range("L31").Value = 4004.6667
$ws_CRP.Range("M31").Value = -1917786.8
$ws_CRP.Range("N31").Value = -4594.6667
$ws_CRP.Range("H34").Value = 1738637
$ws_CRP.Range("I34").Value = 1918081.8
$ws_CRP.Range("J34").Value = 4004.6667
$ws_CRP.Range("K34").Value = 1918081.8
$ws_CRP.Range("L34").Value = 4004.6667
$ws_CRP.Range("M34").Value = -1917879.8
$ws_CRP.Range("N34").Value = -4408.6667
$ws_CRP.Range("H35").Value = 502.16666
$ws_CRP.Range("I35").Value = 502.16666
$ws_CRP.Range("K35").Value = 502.16666
$ws_CRP.Range("M35").Value = -208.16666
$ws_CRP.Range("H39").Value = 500
$ws_CRP.Range("I39").Value = 500
$ws_CRP.Range("K39").Value = 500
$ws_CRP.Range("M39").Value = -109
$ws_CRP.Range("H49").Value = 500
$ws_CRP.Range("I49").Value = 500
$ws_CRP.Range("K49").Value = 500
$ws_CRP.Range("M49").Value = -318
$ws_CRP.Range("H105").Value = 7598.5884
$ws_CRP.Range("I105").Value = 9682.5
$ws_CRP.Range("K105").Value = 9682.5
$ws_CRP.Range("M105").Value = -7935.5
$ws_CRP.Range("H113").Value = 79362.92
$ws_CRP.Range("I113").Value = 2374.2222
$ws_CRP.Range("K113").Value = 2374.2222
$ws_CRP.Range("M113").Value = -204.2222000000002
$ws_CRP.Range("H134").Value = 1498.4546
$ws_CRP.Range("I134").Value = 1379.3334
$ws_CRP.Range("K134").Value = 4138.0002
$ws_CRP.Range("M134").Value = -1603.0002

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H17").Value = 1644.875
$ws_CUL.Range("I17").Value = 552
$ws_CUL.Range("J17").Value = 3466.3333
$ws_CUL.Range("K17").Value = 1656
$ws_CUL.Range("L17").Value = 10398.9999
$ws_CUL.Range("M17").Value = -1487
$ws_CUL.Range("N17").Value = -10736.9999
$ws_CUL.Range("H31").Value = 349
$ws_CUL.Range("I31").Value = 349
$ws_CUL.Range("K31").Value = 1047
$ws_CUL.Range("M31").Value = -759
$ws_CUL.Range("H32").Value = 5050
$ws_CUL.Range("J32").Value = 5050
$ws_CUL.Range("L32").Value = 15150
$ws_CUL.Range("N32").Value = -15716
$ws_CUL.Range("H92").Value = 388.67856
$ws_CUL.Range("J92").Value = 407.57144
$ws_CUL.Range("L92").Value = 1222.71432
$ws_CUL.Range("N92").Value = -3718.71432
$ws_CUL.Range("H122").Value = 621095.75
$ws_CUL.Range("I122").Value = 1241292.9
$ws_CUL.Range("J122").Value = 898.53845
$ws_CUL.Range("K122").Value = 11171636.1
$ws_CUL.Range("L122").Value = 8086.84605
$ws_CUL.Range("M122").Value = -11169186.1
$ws_CUL.Range("N122").Value = -12986.84605
$ws_CUL.Range("H129").Value = 6735537.5
$ws_CUL.Range("I129").Value = 2002921
$ws_CUL.Range("J129").Value = 8555775
$ws_CUL.Range("K129").Value = 6008763
$ws_CUL.Range("L129").Value = 25667325
$ws_CUL.Range("M129").Value = -6003763
$ws_CUL.Range("N129").Value = -25677325

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H113").Value = 2982.0571
$ws_GSM.Range("I113").Value = 2539.52
$ws_GSM.Range("K113").Value = 2539.52
$ws_GSM.Range("M113").Value = -369.52
$ws_GSM.Range("H132").Value = 11098.904
$ws_GSM.Range("I132").Value = 6653.85
$ws_GSM.Range("J132").Value = 100000
$ws_GSM.Range("K132").Value = 19961.55
$ws_GSM.Range("L132").Value = 300000
$ws_GSM.Range("M132").Value = -17431.55
$ws_GSM.Range("N132").Value = -305060

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H46").Value = 3264.8
$ws_LTW.Range("I46").Value = 437.85715
$ws_LTW.Range("K46").Value = 437.85715
$ws_LTW.Range("M46").Value = -249.85715
$ws_LTW.Range("H132").Value = 9739466
$ws_LTW.Range("I132").Value = 9739466
$ws_LTW.Range("K132").Value = 29218398
$ws_LTW.Range("M132").Value = -29215868

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H100").Value = 1230.3334
$ws_WVR.Range("I100").Value = 1230.3334
$ws_WVR.Range("K100").Value = 2460.6668
$ws_WVR.Range("M100").Value = -1919.6668
$ws_WVR.Range("H132").Value = 15154979
$ws_WVR.Range("J132").Value = 3819.6
$ws_WVR.Range("L132").Value = 11458.8
$ws_WVR.Range("N132").Value = -16518.8
